# Applies cryptos list update per commit "Updated cryptos list on Sun Mar 24 16:46:12 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.650.18'
$ws.Range('E2').Value = '  -0.07%  '

$ws.Range('D3').Value = '3.391.77'
$ws.Range('E3').Value = '  -0.98%  '

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '561.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.64%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.32%  '

$ws.Range('E7').Value = '  +0.31%  '

$ws.Range('D8').Value = '3.381.32'
$ws.Range('E8').Value = '  -0.84%  '

$ws.Range('E9').Value = '  +0.05%  '

$ws.Range('E10').Value = '  +2.73%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.633'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.15%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.56'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.82%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000277'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.22%  '

$ws.Range('D15').Value = '3.936.12'
$ws.Range('E15').Value = '  -0.59%  '

$ws.Range('E16').Value = '  -1.05%  '

$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.119'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.14%  '

$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.380.41'
$ws.Range('E18').Value = '  -0.97%  '

$ws.Range('D19').Value = '65.653.95'
$ws.Range('E19').Value = '  +0.22%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.84'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.40%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.998'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.19%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '479.74'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.08%  '

$ws.Range('E23').Value = '  -0.47%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '90.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.82%  '

$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.11'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.46%  '

$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.26%  '

$ws.Range('E27').Value = '  -0.17%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.62'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.08%  '

$ws.Range('E29').Value = '  -2.28%  '

$ws.Range('E30').Value = '  +1.40%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.56'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.75%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '63.80'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.53%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.43'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.33%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '573.54'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.48%  '

$ws.Range('E35').Value = '  -1.27%  '

$ws.Range('E36').Value = '  -0.02%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.62'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.93%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.142'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.06%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.74'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.31%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.374'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.79%  '

$ws.Range('D41').Value = '0.0₃0742'
$ws.Range('E41').Value = '  -3.32%  '

$ws.Range('D42').Value = '3.099.89'
$ws.Range('E42').Value = '  -0.60%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.22%  '

$ws.Range('E44').Value = '  -0.32%  '

$ws.Range('E45').Value = '  -0.48%  '

$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.87%  '

$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.56%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.26%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.58'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.87%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.57'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.14%  '

$ws.Range('E51').Value = '  +0.42%  '
